$wb = $excel.ActiveWorkbook

# --- ALC row 29 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 5000
$ws.Cells.Item(29, 10).Value = 5000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 14).Value = -15562

# --- ALC row 43 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2597.2856
$ws.Cells.Item(43, 9).Value = 1529
$ws.Cells.Item(43, 11).Value = 1529
$ws.Cells.Item(43, 13).Value = -1460

# --- ALC row 58 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 12524.059
$ws.Cells.Item(58, 10).Value = 22579.666
$ws.Cells.Item(58, 12).Value = 67738.99800000001
$ws.Cells.Item(58, 14).Value = -68038.99800000001

# --- ALC row 70 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 4012
$ws.Cells.Item(70, 9).Value = 1609.75
$ws.Cells.Item(70, 10).Value = 5213.125
$ws.Cells.Item(70, 11).Value = 4829.25
$ws.Cells.Item(70, 12).Value = 15639.375
$ws.Cells.Item(70, 13).Value = -4559.25
$ws.Cells.Item(70, 14).Value = -16179.375

# --- ALC row 73 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 4012
$ws.Cells.Item(73, 9).Value = 1609.75
$ws.Cells.Item(73, 10).Value = 5213.125
$ws.Cells.Item(73, 11).Value = 4829.25
$ws.Cells.Item(73, 12).Value = 15639.375
$ws.Cells.Item(73, 13).Value = -3893.25
$ws.Cells.Item(73, 14).Value = -17511.375

# --- ALC row 100 (hunk 5) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 8289.444
$ws.Cells.Item(100, 9).Value = 26549.5
$ws.Cells.Item(100, 11).Value = 26549.5
$ws.Cells.Item(100, 13).Value = -26008.5

# --- ALC row 111 (hunk 6) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1602.2858
$ws.Cells.Item(111, 10).Value = 2438.5
$ws.Cells.Item(111, 12).Value = 7315.5
$ws.Cells.Item(111, 14).Value = -13449.5

# --- ALC row 116 (hunk 7) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 3403
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 3403
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 3403
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(116, 14).Value = -10287

# --- ALC row 131 (hunk 8) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 4898.375
$ws.Cells.Item(131, 10).Value = 8825
$ws.Cells.Item(131, 12).Value = 26475
$ws.Cells.Item(131, 14).Value = -36555

# --- ARM row 6 (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 79320.336
$ws.Cells.Item(6, 9).Value = 92783.8
$ws.Cells.Item(6, 10).Value = 12003
$ws.Cells.Item(6, 11).Value = 92783.8
$ws.Cells.Item(6, 12).Value = 12003
$ws.Cells.Item(6, 13).Value = -92610.8
$ws.Cells.Item(6, 14).Value = -12349

# --- ARM row 19 (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 2501
$ws.Cells.Item(19, 10).Value = 4994
$ws.Cells.Item(19, 12).Value = 4994
$ws.Cells.Item(19, 14).Value = -5452

# --- ARM row 22 (hunk 11) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 1387.8334
$ws.Cells.Item(22, 9).Value = 797.6667
$ws.Cells.Item(22, 10).Value = 1978
$ws.Cells.Item(22, 11).Value = 797.6667
$ws.Cells.Item(22, 12).Value = 1978
$ws.Cells.Item(22, 13).Value = -498.6667
$ws.Cells.Item(22, 14).Value = -2576

# --- ARM row 26 (hunk 12) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 4529.6665
$ws.Cells.Item(26, 9).Value = 4647
$ws.Cells.Item(26, 11).Value = 4647
$ws.Cells.Item(26, 13).Value = -4317

# --- ARM row 27 (hunk 13) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()

# --- ARM row 36 (hunk 14) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 5308
$ws.Cells.Item(36, 9).Value = 5308
$ws.Cells.Item(36, 11).Value = 5308
$ws.Cells.Item(36, 13).Value = -4962

# --- ARM row 61 (hunk 15) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4768.6875
$ws.Cells.Item(61, 9).Value = 4586.6665
$ws.Cells.Item(61, 10).Value = 7499
$ws.Cells.Item(61, 11).Value = 4586.6665
$ws.Cells.Item(61, 12).Value = 7499
$ws.Cells.Item(61, 13).Value = -4374.6665
$ws.Cells.Item(61, 14).Value = -7923

# --- ARM row 136 (hunk 16) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4768.6875
$ws.Cells.Item(136, 9).Value = 4586.6665
$ws.Cells.Item(136, 10).Value = 7499
$ws.Cells.Item(136, 11).Value = 13759.9995
$ws.Cells.Item(136, 12).Value = 22497
$ws.Cells.Item(136, 13).Value = -11209.9995
$ws.Cells.Item(136, 14).Value = -27597

# --- BSM row 7 (hunk 17) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 3807.5
$ws.Cells.Item(7, 9).Value = 2500
$ws.Cells.Item(7, 11).Value = 2500
$ws.Cells.Item(7, 13).Value = -2387

# --- CRP row 12 (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 3627
$ws.Cells.Item(12, 9).Value = 4699
$ws.Cells.Item(12, 10).Value = 2555
$ws.Cells.Item(12, 11).Value = 4699
$ws.Cells.Item(12, 12).Value = 2555
$ws.Cells.Item(12, 13).Value = -4529
$ws.Cells.Item(12, 14).Value = -2895

# --- CRP row 119 (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(119, 8).Value = 53325.668
$ws.Cells.Item(119, 10).Value = 53325.668
$ws.Cells.Item(119, 12).Value = 53325.668
$ws.Cells.Item(119, 14).Value = -63001.668

# --- CUL row 4 (hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 430691.84
$ws.Cells.Item(4, 9).Value = 444
$ws.Cells.Item(4, 11).Value = 1332
$ws.Cells.Item(4, 13).Value = -1220

# --- CUL row 14 (hunk 21) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 75.28570999999999
$ws.Cells.Item(14, 9).Value = 75.28570999999999
$ws.Cells.Item(14, 11).Value = 225.85713
$ws.Cells.Item(14, 13).Value = -52.85712999999998

# --- GSM row 15 (hunk 22) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 9834.75
$ws.Cells.Item(15, 10).Value = 9834.75
$ws.Cells.Item(15, 12).Value = 9834.75
$ws.Cells.Item(15, 14).Value = -10410.75

# --- GSM row 81 (hunk 23) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(81, 8).Value = 9834.75
$ws.Cells.Item(81, 10).Value = 9834.75
$ws.Cells.Item(81, 12).Value = 9834.75
$ws.Cells.Item(81, 14).Value = -11830.75

# --- GSM row 84 (hunk 24) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(84, 8).Value = 9834.75
$ws.Cells.Item(84, 10).Value = 9834.75
$ws.Cells.Item(84, 12).Value = 29504.25
$ws.Cells.Item(84, 14).Value = -39488.25

# --- GSM row 113 (hunk 25) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2266.6428
$ws.Cells.Item(113, 9).Value = 2171.8462
$ws.Cells.Item(113, 10).Value = 3499
$ws.Cells.Item(113, 11).Value = 2171.8462
$ws.Cells.Item(113, 12).Value = 3499
$ws.Cells.Item(113, 13).Value = -1.846199999999953
$ws.Cells.Item(113, 14).Value = -7839

# --- GSM row 132 (hunk 26) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3291.5
$ws.Cells.Item(132, 9).Value = 3816.3333
$ws.Cells.Item(132, 11).Value = 11448.9999
$ws.Cells.Item(132, 13).Value = -8918.999899999999

# --- LTW row 9 (hunk 27) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 9397.4
$ws.Cells.Item(9, 9).Value = 2336
$ws.Cells.Item(9, 10).Value = 19989.5
$ws.Cells.Item(9, 11).Value = 2336
$ws.Cells.Item(9, 12).Value = 19989.5
$ws.Cells.Item(9, 13).Value = -2112
$ws.Cells.Item(9, 14).Value = -20437.5

# --- LTW row 12 (hunk 28) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 11627
$ws.Cells.Item(12, 10).Value = 11627
$ws.Cells.Item(12, 12).Value = 11627
$ws.Cells.Item(12, 14).Value = -11967

# --- LTW row 34 (hunk 29) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 24122.5
$ws.Cells.Item(34, 9).Value = 23490
$ws.Cells.Item(34, 10).Value = 24333.334
$ws.Cells.Item(34, 11).Value = 23490
$ws.Cells.Item(34, 12).Value = 24333.334
$ws.Cells.Item(34, 13).Value = -23318
$ws.Cells.Item(34, 14).Value = -24677.334

# --- LTW row 75 (hunk 30) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(75, 8).Value = 50078
$ws.Cells.Item(75, 9).Value = 50157
$ws.Cells.Item(75, 10).Value = 49999
$ws.Cells.Item(75, 11).Value = 50157
$ws.Cells.Item(75, 12).Value = 49999
$ws.Cells.Item(75, 13).Value = -49221
$ws.Cells.Item(75, 14).Value = -51871

# --- LTW row 78 (hunk 31) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(78, 8).Value = 50078
$ws.Cells.Item(78, 9).Value = 50157
$ws.Cells.Item(78, 10).Value = 49999
$ws.Cells.Item(78, 11).Value = 150471
$ws.Cells.Item(78, 12).Value = 149997
$ws.Cells.Item(78, 13).Value = -145791
$ws.Cells.Item(78, 14).Value = -159357

# --- LTW row 103 (hunk 32) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(103, 8).Value = 13875.25
$ws.Cells.Item(103, 10).Value = 13875.25
$ws.Cells.Item(103, 12).Value = 13875.25
$ws.Cells.Item(103, 14).Value = -16219.25

# --- WVR row 2 (hunk 33) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 15016.5
$ws.Cells.Item(2, 9).Value = 66
$ws.Cells.Item(2, 11).Value = 66
$ws.Cells.Item(2, 13).Value = 46

# --- WVR row 126 (hunk 34) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2545.08
$ws.Cells.Item(126, 9).Value = 1608.375
$ws.Cells.Item(126, 11).Value = 4825.125
$ws.Cells.Item(126, 13).Value = -2355.125

# --- WVR row 132 (hunk 35) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3306.639
$ws.Cells.Item(132, 10).Value = 4142.5713
$ws.Cells.Item(132, 12).Value = 12427.7139
$ws.Cells.Item(132, 14).Value = -17487.7139
